$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    '${created_at}',   # FindText
    $true,             # MatchCase
    $false,            # MatchWholeWord
    $false,            # MatchWildcards
    $false,            # MatchSoundsLike
    $false,            # MatchAllWordForms
    $true,             # Forward
    1,                 # Wrap (wdFindContinue)
    $false,            # Format
    '……………………',       # ReplaceWith
    2                  # Replace (wdReplaceAll)
)
